$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen all 68 columns from 2.16796875 to (closest achievable to) 12.7109375 characters
for ($i = 1; $i -le 68; $i++) {
    $ws.Columns.Item($i).ColumnWidth = 11.833333333333334
}

# Apply the recorded cell value changes (adjacency matrix re-weighting)
$ws.Range("I1").Value = 0
$ws.Range("Q1").Value = 0
$ws.Range("Y1").Value = 0
$ws.Range("Z1").Value = 0
$ws.Range("AW1").Value = 0.6420919103590903
$ws.Range("BB1").Value = 0.61278728278640082
$ws.Range("BI1").Value = 0.76883785900686352
$ws.Range("BJ1").Value = 0
$ws.Range("BP1").Value = 0.50859671214957236
$ws.Range("D2").Value = 0.64243878227429363
$ws.Range("T2").Value = 0
$ws.Range("AK2").Value = 0.72753945484648819
$ws.Range("AS2").Value = 0.68197283347382431
$ws.Range("AU2").Value = 0
$ws.Range("BG2").Value = 0.75704380933583248
$ws.Range("BL2").Value = 0.90629266200181757
$ws.Range("D3").Value = 0
$ws.Range("Y3").Value = 0.98827814049834695
$ws.Range("AA3").Value = 0.95676311399250191
$ws.Range("AM3").Value = 0
$ws.Range("AU3").Value = 0
$ws.Range("AX3").Value = 0.67893814328889768
$ws.Range("BB3").Value = 0.58804587552326049
$ws.Range("BK3").Value = 0.5752196747711309
$ws.Range("BN3").Value = 0.70841508452175717
$ws.Range("B4").Value = 0.63537838367930322
$ws.Range("C4").Value = 0
$ws.Range("O4").Value = 0.69921353163424682
$ws.Range("R4").Value = 0.99124050818327125
$ws.Range("AA4").Value = 0.9563617453344726
$ws.Range("AL4").Value = 0
$ws.Range("AA5").Value = 0.78496263303642033
$ws.Range("AM5").Value = 0
$ws.Range("BA5").Value = 0.61493076884793663
$ws.Range("BE5").Value = 0
$ws.Range("BG5").Value = 0
$ws.Range("H6").Value = 0.61661260248735084
$ws.Range("N6").Value = 0
$ws.Range("P6").Value = 0.61966163156690657
$ws.Range("AI6").Value = 0
$ws.Range("AP6").Value = 0.80915034199157043
$ws.Range("BA6").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0.78265243085276437
$ws.Range("R7").Value = 0.77676947319056056
$ws.Range("AE7").Value = 0
$ws.Range("AK7").Value = 0
$ws.Range("AQ7").Value = 0
$ws.Range("AV7").Value = 0
$ws.Range("BG7").Value = 0.62107229935074892
$ws.Range("F8").Value = 0.69247908964704641
$ws.Range("G8").Value = 0
$ws.Range("Q8").Value = 0.60307648233710298
$ws.Range("AA8").Value = 0.89701277270397783
$ws.Range("AN8").Value = 0
$ws.Range("AP8").Value = 0.93864992630987332
$ws.Range("BD8").Value = 0
$ws.Range("BL8").Value = 0.98894100599759271
$ws.Range("A9").Value = 0
$ws.Range("G9").Value = 0.82614197959063351
$ws.Range("Y9").Value = 0
$ws.Range("AE9").Value = 0.99641414074189782
$ws.Range("AG9").Value = 0.77535356548840073
$ws.Range("AO9").Value = 0.58017636755472402
$ws.Range("AX9").Value = 0
$ws.Range("K10").Value = 0.83048781503696067
$ws.Range("N10").Value = 0.53190224253608331
$ws.Range("P10").Value = 0
$ws.Range("BE10").Value = 0.83250885266459695
$ws.Range("BL10").Value = 0
$ws.Range("J11").Value = 0.93459226183390798
$ws.Range("M11").Value = 0
$ws.Range("AB11").Value = 0
$ws.Range("AQ11").Value = 0.625891958838662
$ws.Range("AS11").Value = 0.71843609649376905
$ws.Range("AX11").Value = 0.8580495581362283
$ws.Range("BA11").Value = 0.9237682868317284
$ws.Range("M12").Value = 0
$ws.Range("P12").Value = 0.83988054207341967
$ws.Range("Y12").Value = 0.74701662664733992
$ws.Range("AF12").Value = 0.79409793231113834
$ws.Range("AG12").Value = 0
$ws.Range("AN12").Value = 0.75072792625412654
$ws.Range("AV12").Value = 0.72662239230744918
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("R13").Value = 0.86940033184227816
$ws.Range("T13").Value = 0
$ws.Range("V13").Value = 0.934885263985614
$ws.Range("X13").Value = 0.76597302698463654
$ws.Range("AH13").Value = 0.515274783610006
$ws.Range("AR13").Value = 0
$ws.Range("AY13").Value = 0
$ws.Range("BG13").Value = 0.79160632050094937
$ws.Range("F14").Value = 0
$ws.Range("J14").Value = 0.93256332555256061
$ws.Range("P14").Value = 0
$ws.Range("U14").Value = 0
$ws.Range("AI14").Value = 0.81439913200292602
$ws.Range("BD14").Value = 0.67722449560728593
$ws.Range("BK14").Value = 0.89502085826048572
$ws.Range("D15").Value = 0.6181499168686645
$ws.Range("P15").Value = 0.96272025594045529
$ws.Range("AE15").Value = 0
$ws.Range("AG15").Value = 0
$ws.Range("AO15").Value = 0.87841503401033705
$ws.Range("BD15").Value = 0
$ws.Range("BE15").Value = 0.72585186337741658
$ws.Range("BJ15").Value = 0
$ws.Range("F16").Value = 0.99434153047568397
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0.57066265561865737
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 0.69116569760842828
$ws.Range("Q16").Value = 0
$ws.Range("S16").Value = 0
$ws.Range("W16").Value = 0
$ws.Range("AB16").Value = 0
$ws.Range("AD16").Value = 0
$ws.Range("AH16").Value = 0.88847447232631716
$ws.Range("AL16").Value = 0
$ws.Range("AU16").Value = 0.731255128966386
$ws.Range("AZ16").Value = 0
$ws.Range("BB16").Value = 0.57729591727055318
$ws.Range("A17").Value = 0
$ws.Range("H17").Value = 0.93732748956306389
$ws.Range("P17").Value = 0
$ws.Range("R17").Value = 0
$ws.Range("S17").Value = 0
$ws.Range("AO17").Value = 0
$ws.Range("BD17").Value = 0.94133223070222916
$ws.Range("BH17").Value = 0.90250293315166674
$ws.Range("D18").Value = 0.57177197347504438
$ws.Range("G18").Value = 0.93939028850615225
$ws.Range("M18").Value = 0.85374129400172316
$ws.Range("Q18").Value = 0
$ws.Range("T18").Value = 0
$ws.Range("U18").Value = 0.79988425663442664
$ws.Range("AW18").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 0
$ws.Range("U19").Value = 0
$ws.Range("AL19").Value = 0
$ws.Range("AX19").Value = 0.61851860686659554
$ws.Range("BF19").Value = 0.86938179828319095
$ws.Range("BO19").Value = 0.88269159305700828
$ws.Range("B20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("R20").Value = 0
$ws.Range("V20").Value = 0
$ws.Range("AA20").Value = 0.94827605383369373
$ws.Range("AB20").Value = 0
$ws.Range("AF20").Value = 0.51851043083645698
$ws.Range("AX20").Value = 0
$ws.Range("BI20").Value = 0.52769239759498077
$ws.Range("BK20").Value = 0
$ws.Range("BO20").Value = 0.72063045873370579
$ws.Range("N21").Value = 0
$ws.Range("R21").Value = 0.52639671441848734
$ws.Range("S21").Value = 0
$ws.Range("V21").Value = 0.71190555195097127
$ws.Range("W21").Value = 0
$ws.Range("AD21").Value = 0.80806103979724597
$ws.Range("AS21").Value = 0.5708869341057512
$ws.Range("AW21").Value = 0
$ws.Range("BI21").Value = 0.57914681549874725
$ws.Range("M22").Value = 0.56238385848947403
$ws.Range("T22").Value = 0
$ws.Range("U22").Value = 0.78758714264860585
$ws.Range("W22").Value = 0.88535043114292122
$ws.Range("X22").Value = 0.68769893715605557
$ws.Range("AB22").Value = 0
$ws.Range("AL22").Value = 0
$ws.Range("BI22").Value = 0.90264277447611052
$ws.Range("BM22").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("U23").Value = 0
$ws.Range("V23").Value = 0.95002322457079957
$ws.Range("X23").Value = 0
$ws.Range("Z23").Value = 0.65268440521854965
$ws.Range("AT23").Value = 0.94683279029440093
$ws.Range("BK23").Value = 0.73582731427577941
$ws.Range("BN23").Value = 0
$ws.Range("M24").Value = 0.88978598926315966
$ws.Range("V24").Value = 0.73037739070690311
$ws.Range("W24").Value = 0
$ws.Range("Z24").Value = 0
$ws.Range("AV24").Value = 0
$ws.Range("BN24").Value = 0
$ws.Range("A25").Value = 0
$ws.Range("C25").Value = 0.63583089334680354
$ws.Range("I25").Value = 0
$ws.Range("L25").Value = 0.58281878711298463
$ws.Range("Z25").Value = 0
$ws.Range("AM25").Value = 0
$ws.Range("BN25").Value = 0.67937453422822558
$ws.Range("A26").Value = 0
$ws.Range("W26").Value = 0.95354806876194442
$ws.Range("X26").Value = 0
$ws.Range("Y26").Value = 0
$ws.Range("AC26").Value = 0.79034906009574946
$ws.Range("BM26").Value = 0
$ws.Range("BO26").Value = 0.52636312505013239
$ws.Range("C27").Value = 0.63474991947209292
$ws.Range("D27").Value = 0.56382890435780753
$ws.Range("E27").Value = 0.98206569606794747
$ws.Range("H27").Value = 0.99079801676383794
$ws.Range("T27").Value = 0.76513835336430813
$ws.Range("AB27").Value = 0.59870721397640714
$ws.Range("AC27").Value = 0
$ws.Range("AQ27").Value = 0.76127168325092631
$ws.Range("AU27").Value = 0
$ws.Range("BJ27").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("P28").Value = 0
$ws.Range("T28").Value = 0
$ws.Range("V28").Value = 0
$ws.Range("AA28").Value = 0.85284958888532181
$ws.Range("AO28").Value = 0.68390147084585218
$ws.Range("AY28").Value = 0.55359536118111574
$ws.Range("AZ28").Value = 0
$ws.Range("BB28").Value = 0
$ws.Range("Z29").Value = 0.59376427527978415
$ws.Range("AA29").Value = 0
$ws.Range("AE29").Value = 0
$ws.Range("BD29").Value = 0.99417741078455779
$ws.Range("BF29").Value = 0.63847787478756901
$ws.Range("BK29").Value = 0
$ws.Range("BO29").Value = 0.59672601466091879
$ws.Range("P30").Value = 0
$ws.Range("U30").Value = 0.89735322424683051
$ws.Range("AE30").Value = 0.98342541977692366
$ws.Range("AI30").Value = 0
$ws.Range("BF30").Value = 0.55973331201012644
$ws.Range("G31").Value = 0
$ws.Range("I31").Value = 0.93377463020275497
$ws.Range("O31").Value = 0
$ws.Range("AC31").Value = 0
$ws.Range("AD31").Value = 0.8148849410999075
$ws.Range("AF31").Value = 0.73076308806403456
$ws.Range("AG31").Value = 0
$ws.Range("AN31").Value = 0.79520073956512971
$ws.Range("AP31").Value = 0.94847474043333069
$ws.Range("AR31").Value = 0
$ws.Range("BB31").Value = 0
$ws.Range("L32").Value = 0.51274243966165023
$ws.Range("T32").Value = 0.63681178827866436
$ws.Range("AE32").Value = 0.93851088690433437
$ws.Range("AR32").Value = 0.90398502543581749
$ws.Range("AT32").Value = 0
$ws.Range("I33").Value = 0.63406925665564295
$ws.Range("L33").Value = 0
$ws.Range("O33").Value = 0
$ws.Range("AE33").Value = 0
$ws.Range("AH33").Value = 0
$ws.Range("AI33").Value = 0.93099109109474409
$ws.Range("BA33").Value = 0.86473109788015223
$ws.Range("M34").Value = 0.69597876362001498
$ws.Range("P34").Value = 0.84905811272138232
$ws.Range("AG34").Value = 0
$ws.Range("AM34").Value = 0.63247478001399793
$ws.Range("AT34").Value = 0.82908404156231019
$ws.Range("BC34").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("N35").Value = 0.64685705345582645
$ws.Range("AD35").Value = 0
$ws.Range("AG35").Value = 0.93512639783096874
$ws.Range("AJ35").Value = 0.89930495113074482
$ws.Range("AK35").Value = 0.82421796846616124
$ws.Range("AR35").Value = 0
$ws.Range("AT35").Value = 0.60153189490931924
$ws.Range("AI36").Value = 0.71096952784148926
$ws.Range("AL36").Value = 0.69609695686198281
$ws.Range("AM36").Value = 0.59631637682634819
$ws.Range("AP36").Value = 0
$ws.Range("BL36").Value = 0.78387988707644674
$ws.Range("B37").Value = 0.60542810443054806
$ws.Range("G37").Value = 0
$ws.Range("AI37").Value = 0.92521020644400243
$ws.Range("AL37").Value = 0.68083773531089209
$ws.Range("AO37").Value = 0.64163502770043057
$ws.Range("AY37").Value = 0.55397086998553546
$ws.Range("BB37").Value = 0
$ws.Range("D38").Value = 0
$ws.Range("P38").Value = 0
$ws.Range("S38").Value = 0
$ws.Range("V38").Value = 0
$ws.Range("AJ38").Value = 0.61113413117717474
$ws.Range("AK38").Value = 0.95595528442281386
$ws.Range("AM38").Value = 0.53957641777461274
$ws.Range("BC38").Value = 0
$ws.Range("C39").Value = 0
$ws.Range("E39").Value = 0
$ws.Range("Y39").Value = 0
$ws.Range("AH39").Value = 0.90367782859291435
$ws.Range("AJ39").Value = 0.6227499458987884
$ws.Range("AL39").Value = 0.68153274621581283
$ws.Range("H40").Value = 0
$ws.Range("L40").Value = 0.66511663436851298
$ws.Range("AE40").Value = 0.69469273012309962
$ws.Range("BC40").Value = 0
$ws.Range("BF40").Value = 0
$ws.Range("BI40").Value = 0
$ws.Range("BN40").Value = 0
$ws.Range("I41").Value = 0.75577886138754835
$ws.Range("O41").Value = 0.92602854960673064
$ws.Range("Q41").Value = 0
$ws.Range("AB41").Value = 0.56318293491072358
$ws.Range("AK41").Value = 0.93482621032537438
$ws.Range("F42").Value = 0.93504637356407261
$ws.Range("H42").Value = 0.74244861907274662
$ws.Range("AE42").Value = 0.85432805501148423
$ws.Range("AJ42").Value = 0
$ws.Range("AR42").Value = 0
$ws.Range("AS42").Value = 0
$ws.Range("AU42").Value = 0.7323333287898417
$ws.Range("AX42").Value = 0.54154842101952694
$ws.Range("BD42").Value = 0
$ws.Range("BN42").Value = 0.67989051822041002
$ws.Range("G43").Value = 0
$ws.Range("K43").Value = 0.82175026593008416
$ws.Range("AA43").Value = 0.82965778910465582
$ws.Range("AR43").Value = 0.84515735830832139
$ws.Range("AS43").Value = 0
$ws.Range("AZ43").Value = 0.7711888288520774
$ws.Range("BJ43").Value = 0.80400321486404414
$ws.Range("BM43").Value = 0.64065150216380107
$ws.Range("M44").Value = 0
$ws.Range("AE44").Value = 0
$ws.Range("AF44").Value = 0.97041584285706839
$ws.Range("AI44").Value = 0
$ws.Range("AP44").Value = 0
$ws.Range("AQ44").Value = 0.65940112112115457
$ws.Range("AS44").Value = 0
$ws.Range("BA44").Value = 0
$ws.Range("BL44").Value = 0
$ws.Range("BO44").Value = 0
$ws.Range("B45").Value = 0.81779974890632201
$ws.Range("K45").Value = 0.84867351054303186
$ws.Range("U45").Value = 0.86426847102747717
$ws.Range("AP45").Value = 0
$ws.Range("AQ45").Value = 0
$ws.Range("AR45").Value = 0
$ws.Range("AU45").Value = 0
$ws.Range("W46").Value = 0.9138114595129132
$ws.Range("AF46").Value = 0
$ws.Range("AH46").Value = 0.84248493313462802
$ws.Range("AI46").Value = 0.9721754361411894
$ws.Range("AU46").Value = 0
$ws.Range("AV46").Value = 0
$ws.Range("B47").Value = 0
$ws.Range("C47").Value = 0
$ws.Range("P47").Value = 0.84952161440462493
$ws.Range("AA47").Value = 0
$ws.Range("AP47").Value = 0.57293033272781058
$ws.Range("AS47").Value = 0
$ws.Range("AT47").Value = 0
$ws.Range("BL47").Value = 0.57798806065929509
$ws.Range("G48").Value = 0
$ws.Range("L48").Value = 0.83967797284742418
$ws.Range("X48").Value = 0
$ws.Range("AT48").Value = 0
$ws.Range("AW48").Value = 0.63025874371372848
$ws.Range("A49").Value = 0.6759582412006081
$ws.Range("R49").Value = 0
$ws.Range("U49").Value = 0
$ws.Range("AV49").Value = 0.69692563973706512
$ws.Range("AY49").Value = 0.98602939259373812
$ws.Range("BM49").Value = 0
$ws.Range("C50").Value = 0.9126217709345581
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0.71228308607711632
$ws.Range("S50").Value = 0.80826682555525875
$ws.Range("T50").Value = 0
$ws.Range("AP50").Value = 0.96771476716851557
$ws.Range("AY50").Value = 0
$ws.Range("BD50").Value = 0
$ws.Range("M51").Value = 0
$ws.Range("AB51").Value = 0.53235793668444864
$ws.Range("AK51").Value = 0.83808997383320938
$ws.Range("AW51").Value = 0.59854587551923322
$ws.Range("AX51").Value = 0
$ws.Range("AZ51").Value = 0.77551873539298044
$ws.Range("BA51").Value = 0.89015475180454695
$ws.Range("BD51").Value = 0.68174817825292644
$ws.Range("BJ51").Value = 0
$ws.Range("P52").Value = 0
$ws.Range("AB52").Value = 0
$ws.Range("AQ52").Value = 0.5932303843190414
$ws.Range("AY52").Value = 0.67372889648186196
$ws.Range("BP52").Value = 0.61412615287226457
$ws.Range("E53").Value = 0.68848822567416224
$ws.Range("F53").Value = 0
$ws.Range("K53").Value = 0.9017996958458796
$ws.Range("AG53").Value = 0.91752327024113645
$ws.Range("AR53").Value = 0
$ws.Range("AY53").Value = 0.69949404699504214
$ws.Range("BF53").Value = 0
$ws.Range("BO53").Value = 0.53404018897223537
$ws.Range("A54").Value = 0.613735526658579
$ws.Range("C54").Value = 0.60096268347368997
$ws.Range("P54").Value = 0.76898866745847738
$ws.Range("AB54").Value = 0
$ws.Range("AE54").Value = 0
$ws.Range("AK54").Value = 0
$ws.Range("BD54").Value = 0
$ws.Range("AH55").Value = 0
$ws.Range("AL55").Value = 0
$ws.Range("AN55").Value = 0
$ws.Range("BI55").Value = 0.82125914621751073
$ws.Range("BK55").Value = 0.60631569803994556
$ws.Range("BL55").Value = 0
$ws.Range("H56").Value = 0
$ws.Range("N56").Value = 0.53032586904261503
$ws.Range("O56").Value = 0
$ws.Range("Q56").Value = 0.83817592228923365
$ws.Range("AC56").Value = 0.91465531777311626
$ws.Range("AP56").Value = 0
$ws.Range("AX56").Value = 0
$ws.Range("AY56").Value = 0.55239323663048667
$ws.Range("BB56").Value = 0
$ws.Range("BF56").Value = 0.66129536881574325
$ws.Range("BG56").Value = 0
$ws.Range("BK56").Value = 0
$ws.Range("E57").Value = 0
$ws.Range("J57").Value = 0.83600180773863941
$ws.Range("O57").Value = 0.75265914331660633
$ws.Range("BF57").Value = 0.91820669686115242
$ws.Range("BM57").Value = 0.54699328675351944
$ws.Range("BO57").Value = 0
$ws.Range("BP57").Value = 0
$ws.Range("S58").Value = 0.68799019014022578
$ws.Range("AC58").Value = 0.97109735972312072
$ws.Range("AD58").Value = 0.75084428501549416
$ws.Range("AN58").Value = 0
$ws.Range("BA58").Value = 0
$ws.Range("BD58").Value = 0.82385910023990894
$ws.Range("BE58").Value = 0.81643906051908943
$ws.Range("BH58").Value = 0
$ws.Range("BK58").Value = 0.63219084501569878
$ws.Range("B59").Value = 0.56116522806002833
$ws.Range("E59").Value = 0
$ws.Range("G59").Value = 0.83934135952786848
$ws.Range("M59").Value = 0.96792955025418748
$ws.Range("BD59").Value = 0
$ws.Range("BH59").Value = 0
$ws.Range("BI59").Value = 0.69404543668559293
$ws.Range("Q60").Value = 0.5356662939343988
$ws.Range("BF60").Value = 0
$ws.Range("BG60").Value = 0
$ws.Range("BJ60").Value = 0
$ws.Range("BO60").Value = 0.61822738411318734
$ws.Range("A61").Value = 0.67904799847658692
$ws.Range("T61").Value = 0.5190774402993692
$ws.Range("U61").Value = 0.61453340484913122
$ws.Range("V61").Value = 0.92308014836763674
$ws.Range("AN61").Value = 0
$ws.Range("BC61").Value = 0.54239217236607828
$ws.Range("BG61").Value = 0.88904772460884285
$ws.Range("BJ61").Value = 0.79220074209631042
$ws.Range("BO61").Value = 0
$ws.Range("A62").Value = 0
$ws.Range("O62").Value = 0
$ws.Range("AA62").Value = 0
$ws.Range("AQ62").Value = 0.62288523367402615
$ws.Range("AY62").Value = 0
$ws.Range("BH62").Value = 0
$ws.Range("BI62").Value = 0.88743736724617173
$ws.Range("BK62").Value = 0.89927843738899738
$ws.Range("C63").Value = 0.77718134286748342
$ws.Range("N63").Value = 0.78204011182485789
$ws.Range("T63").Value = 0
$ws.Range("W63").Value = 0.78497327705881237
$ws.Range("AC63").Value = 0
$ws.Range("BC63").Value = 0.8111245283993983
$ws.Range("BD63").Value = 0
$ws.Range("BF63").Value = 0.85383546309846869
$ws.Range("BJ63").Value = 0.61285594222998185
$ws.Range("BM63").Value = 0
$ws.Range("B64").Value = 0.76006952572508113
$ws.Range("H64").Value = 0.74557811265807983
$ws.Range("J64").Value = 0
$ws.Range("AJ64").Value = 0.60004162036399245
$ws.Range("AR64").Value = 0
$ws.Range("AU64").Value = 0.80142222521930828
$ws.Range("BC64").Value = 0
$ws.Range("V65").Value = 0
$ws.Range("Z65").Value = 0
$ws.Range("AQ65").Value = 0.7353119599092599
$ws.Range("AW65").Value = 0
$ws.Range("BE65").Value = 0.78640530335439873
$ws.Range("BK65").Value = 0
$ws.Range("BN65").Value = 0.78810765208627964
$ws.Range("C66").Value = 0.61446819235273087
$ws.Range("W66").Value = 0
$ws.Range("X66").Value = 0
$ws.Range("Y66").Value = 0.88601589755283894
$ws.Range("AN66").Value = 0
$ws.Range("AP66").Value = 0.67472446563247435
$ws.Range("BM66").Value = 0.62888981332384231
$ws.Range("BP66").Value = 0
$ws.Range("S67").Value = 0.59645007178622722
$ws.Range("T67").Value = 0.57586395135496904
$ws.Range("Z67").Value = 0.76279126365547445
$ws.Range("AC67").Value = 0.99804711834209059
$ws.Range("AR67").Value = 0
$ws.Range("BA67").Value = 0.75178612163005853
$ws.Range("BE67").Value = 0
$ws.Range("BH67").Value = 0.52073429562560225
$ws.Range("BI67").Value = 0
$ws.Range("BP67").Value = 0.77692949726630967
$ws.Range("A68").Value = 0.82485008018293737
$ws.Range("AZ68").Value = 0.50827692275953007
$ws.Range("BE68").Value = 0
$ws.Range("BN68").Value = 0
$ws.Range("BO68").Value = 0.53400365613997636
